# "Generate Report for Archive"
#
# The localization status report is being regenerated: every cell that
# still shows the old "Ready for handoff" status is now "In Translation"
# (Overview!E:F for zh-cn/de-de, and the Status column on each locale
# sheet). Excel then auto-sized the Status/locale columns to the new
# (shorter) text, so their stored widths shrink as well.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: compare with the literal string on the LEFT of -eq. Some
            # cells hold booleans ("True"/"False" text that Value2 surfaces
            # as System.Boolean); PowerShell's -eq coerces the right-hand
            # side to the left operand's type, so "$boolCell -eq $string"
            # would wrongly coerce the string to a (truthy) bool and match.
            # Putting the known string first keeps the comparison a true
            # string comparison regardless of the cell's underlying type.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Re-fit the columns that held the status text to their new (narrower)
# width, matching what Excel's column AutoFit produces for the shorter
# "In Translation" label.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns("E:F").ColumnWidth = 12.54

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns("C:C").ColumnWidth = 12.54

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns("C:C").ColumnWidth = 12.54
